$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: new activity "3. iterace - diagramy komponent" with 1 hour logged
$ws.Range("A29").Value = "3. iterace - diagramy komponent"
$ws.Range("B29").Value = 1

# Row 30: new activity "3. iterace - návrhové třídy" with 1 hour logged
$ws.Range("A30").Value = "3. iterace - návrhové třídy"
$ws.Range("B30").Value = 1

# Match the formatting used by the other activity name cells above (A27:A28)
$ws.Range("A27:A28").Copy()
$ws.Range("A29:A30").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The active selection moved on to the next empty row after data entry
$ws.Range("B31").Select()
